# Generate Report for Archive
#
# 1. Status text: "Ready for handoff" -> "In Translation" everywhere it
#    appears (Overview!E2:F3 and the "Status" column of the per-locale
#    sheets).
# 2. Narrow the "Status" column(s) from ~17.2 chars to ~13.4 chars on all
#    three sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: zh-cn / de-de status columns (E, F) -----------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"

# ColumnWidth is expressed in characters and gets snapped to the sheet's
# pixel grid, same as real Excel - feed it the inverse of that grid
# mapping so the stored width lands as close as possible to 13.4101845877511.
$wsOverview.Range("E1").ColumnWidth = 12.576851254417766
$wsOverview.Range("F1").ColumnWidth = 12.576851254417766

# --- zh-cn sheet: Status column (C) ---------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("C3").Value = "In Translation"
$wsZhCn.Range("C1").ColumnWidth = 12.576851254417766

# --- de-de sheet: Status column (C) ---------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("C3").Value = "In Translation"
$wsDeDe.Range("C1").ColumnWidth = 12.576851254417766
